# Applies the cryptos.xlsx price/volume/coin-order update described in the commit
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '20.018.89'
$ws.Range('E2').Value = '  -7.38%  '
$ws.Range('D3').Value = '1.420.19'
$ws.Range('E3').Value = '  -7.40%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.000'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.26%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.9960'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.55%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '274.08'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -5.45%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3685'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -5.64%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3133'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.68%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '39.95'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -6.96%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.044'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -2.53%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.06509'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -9.21%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.001'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.23%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.496'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -4.36%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '17.77'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.14%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.204'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -5.83%  '
$ws.Range('B16').Value = 'WrappedEther'
$ws.Range('C16').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D16').Value = '1.417.87'
$ws.Range('E16').Value = '  -7.94%  '
$ws.Range('B17').Value = 'ShibaInu'
$ws.Range('C17').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001018'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -6.60%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.05702'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -13.99%  '
$ws.Range('B19').Value = 'Litecoin'
$ws.Range('C19').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '71.33'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -15.27%  '
$ws.Range('B20').Value = 'Dai'
$ws.Range('C20').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.9967'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.49%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.624'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -8.12%  '
$ws.Range('E22').Value = '  -4.10%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '11.17'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +4.11%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.272'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -4.88%  '
$ws.Range('D25').Value = '20.089.15'
$ws.Range('E25').Value = '  -7.09%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.276'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -3.37%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '135.75'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -10.64%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '17.20'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -6.69%  '
$ws.Range('D29').Value = '1.579.54'
$ws.Range('E29').Value = '  -7.67%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '110.36'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -5.79%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.006'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -18.26%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '5.387'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -10.61%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.8385'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -11.54%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.07725'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -3.84%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '8.446'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.41%  '
$ws.Range('B36').Value = 'WEMIXTOKEN'
$ws.Range('C36').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.477'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -0.67%  '
$ws.Range('B37').Value = 'Hedera'
$ws.Range('C37').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.05919'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +0.49%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.893'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -5.42%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.9960'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -0.50%  '
$ws.Range('B40').Value = 'Aptos'
$ws.Range('C40').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '10.71'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -4.61%  '
$ws.Range('B41').Value = 'VeChain'
$ws.Range('C41').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.02075'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -5.94%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.1919'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -5.67%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.098'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -7.22%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.5317'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -7.73%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '12.37'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -6.00%  '
$ws.Range('E46').Value = '  -5.25%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.5179'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -6.67%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '114.63'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.10%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.772'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -6.01%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.046'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -9.45%  '
$ws.Range('B51').Value = 'Cronos'
$ws.Range('C51').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.06237'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -7.01%  '
